$wb = $excel.ActiveWorkbook

# Sheets 1-5 (PagoSinTarjetaAsociada, PagoAfiliadoDebitoAuto, PagoConValidacionHistorial,
# PagoConTarjetaAsociada, PagoAsociandoTarjeta) all share the same test data row that is
# being refreshed: numeroUsuario, placa and the last date (fecha/vigencia) column.
# PagoConTarjetaAsociada (sheet4 / Tabla36) only has columns A:E (vigencia is column E),
# the rest use columns A:G (vigencia is column G).
$sheetColumns = @{
    "PagoSinTarjetaAsociada"      = "G"
    "PagoAfiliadoDebitoAuto"      = "G"
    "PagoConValidacionHistorial"  = "G"
    "PagoConTarjetaAsociada"      = "E"
    "PagoAsociandoTarjeta"        = "G"
}

foreach ($name in $sheetColumns.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A2").Value = "'72636759"
    $ws.Range("C2").Value = "'ZIJ-583"
    $dateCol = $sheetColumns[$name]
    $ws.Range($dateCol + "2").Value = "'12/11/2025"
}

# Update the selection anchors so they match what the author left selected in Excel.
# (Activate() + Select() are used per-sheet so each sheet keeps its own cached
# selection; the original active tab is restored afterwards.)
$wb.Worksheets.Item("PagoAfiliadoDebitoAuto").Activate()
$wb.Worksheets.Item("PagoAfiliadoDebitoAuto").Range("G2").Select()

$wb.Worksheets.Item("PagoConValidacionHistorial").Activate()
$wb.Worksheets.Item("PagoConValidacionHistorial").Range("G2").Select()

$wb.Worksheets.Item("PagoInfoNoDisponible").Activate()
$wb.Worksheets.Item("PagoInfoNoDisponible").Range("C5").Select()

# Restore the originally active sheet (5th tab, PagoAsociandoTarjeta).
$wb.Worksheets.Item("PagoAsociandoTarjeta").Activate()
